$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 660111.6
$ws.Range("I15").Value = 660111.6
$ws.Range("K15").Value = 1980334.8
$ws.Range("M15").Value = -1980165.8
$ws.Range("H26").Value = 4000
$ws.Range("J26").Value = 4000
$ws.Range("L26").Value = 4000
$ws.Range("N26").Value = -4688
$ws.Range("H53").Value = 3580
$ws.Range("J53").Value = 350
$ws.Range("L53").Value = 350
$ws.Range("N53").Value = -1624
$ws.Range("H74").Value = 4664
$ws.Range("I74").Value = 4608
$ws.Range("K74").Value = 4608
$ws.Range("M74").Value = -3672
$ws.Range("H77").Value = 4664
$ws.Range("I77").Value = 4608
$ws.Range("K77").Value = 23040
$ws.Range("M77").Value = -18360
$ws.Range("H100").Value = 17966364
$ws.Range("I100").Value = 31376962
$ws.Range("K100").Value = 31376962
$ws.Range("M100").Value = -31376421
$ws.Range("H113").Value = 19184.857
$ws.Range("I113").Value = 18258.8
$ws.Range("K113").Value = 18258.8
$ws.Range("M113").Value = -15004.8
$ws.Range("H116").Value = 8782040
$ws.Range("I116").Value = 19021384
$ws.Range("J116").Value = 5460.4287
$ws.Range("K116").Value = 19021384
$ws.Range("L116").Value = 5460.4287
$ws.Range("M116").Value = -19017942
$ws.Range("N116").Value = -12344.4287
$ws.Range("H137").Value = 752515.3
$ws.Range("I137").Value = 1064733.2
$ws.Range("J137").Value = 3192.1
$ws.Range("K137").Value = 3194199.6
$ws.Range("L137").Value = 9576.299999999999
$ws.Range("M137").Value = -3191649.6
$ws.Range("N137").Value = -14676.3
$ws.Range("H138").Value = 163175.84
$ws.Range("I138").Value = 694426
$ws.Range("J138").Value = 5460.953
$ws.Range("K138").Value = 2083278
$ws.Range("L138").Value = 16382.859
$ws.Range("M138").Value = -2078138
$ws.Range("N138").Value = -26662.859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 97063.23
$ws.Range("I45").Value = 123199.47
$ws.Range("J45").Value = 8200
$ws.Range("K45").Value = 123199.47
$ws.Range("L45").Value = 8200
$ws.Range("M45").Value = -122822.47
$ws.Range("N45").Value = -8954
$ws.Range("H63").Value = 2864.7144
$ws.Range("I63").Value = 2864.7144
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2864.7144
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2178.7144
$ws.Range("H66").Value = 2864.7144
$ws.Range("I66").Value = 2864.7144
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14323.572
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10891.572
$ws.Range("H80").Value = 84999.86
$ws.Range("J80").Value = 84999.86
$ws.Range("L80").Value = 84999.86
$ws.Range("N80").Value = -86995.86
$ws.Range("H83").Value = 84999.86
$ws.Range("J83").Value = 84999.86
$ws.Range("L83").Value = 254999.58
$ws.Range("N83").Value = -264983.58
$ws.Range("H97").Value = 11770871
$ws.Range("I97").Value = 6550.875
$ws.Range("K97").Value = 6550.875
$ws.Range("M97").Value = -6054.875
$ws.Range("H122").Value = 1586066.8
$ws.Range("I122").Value = 7566.077
$ws.Range("K122").Value = 22698.231
$ws.Range("M122").Value = -20248.231
$ws.Range("H132").Value = 4038.1
$ws.Range("I132").Value = 3105.3333
$ws.Range("K132").Value = 9315.999899999999
$ws.Range("M132").Value = -6785.999899999999
$ws.Range("H134").Value = 237497
$ws.Range("J134").Value = 237497
$ws.Range("L134").Value = 237497
$ws.Range("N134").Value = -247637
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1665.1
$ws.Range("I20").Value = 1096.6428
$ws.Range("K20").Value = 1096.6428
$ws.Range("M20").Value = -849.6428000000001
$ws.Range("H134").Value = 2280.2927
$ws.Range("I134").Value = 1710.6875
$ws.Range("J134").Value = 4305.5557
$ws.Range("K134").Value = 5132.0625
$ws.Range("L134").Value = 12916.6671
$ws.Range("M134").Value = -2597.0625
$ws.Range("N134").Value = -17986.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3390.8333
$ws.Range("I31").Value = 2513.3333
$ws.Range("J31").Value = 4707.0835
$ws.Range("K31").Value = 2513.3333
$ws.Range("L31").Value = 4707.0835
$ws.Range("M31").Value = -2218.3333
$ws.Range("N31").Value = -5297.0835
$ws.Range("H34").Value = 3390.8333
$ws.Range("I34").Value = 2513.3333
$ws.Range("J34").Value = 4707.0835
$ws.Range("K34").Value = 2513.3333
$ws.Range("L34").Value = 4707.0835
$ws.Range("M34").Value = -2311.3333
$ws.Range("N34").Value = -5111.0835
$ws.Range("H105").Value = 194815.64
$ws.Range("I105").Value = 213797.2
$ws.Range("K105").Value = 213797.2
$ws.Range("M105").Value = -212050.2
$ws.Range("H132").Value = 10136.363
$ws.Range("I132").Value = 11900
$ws.Range("K132").Value = 35700
$ws.Range("M132").Value = -33170
$ws.Range("H134").Value = 2548.1924
$ws.Range("I134").Value = 2493.8262
$ws.Range("K134").Value = 7481.4786
$ws.Range("M134").Value = -4946.4786
$ws.Range("H141").Value = 369060.34
$ws.Range("J141").Value = 438134.7
$ws.Range("L141").Value = 438134.7
$ws.Range("N141").Value = -448494.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 6000
$ws.Range("I120").Value = 6000
$ws.Range("K120").Value = 18000
$ws.Range("M120").Value = -13162
$ws.Range("H122").Value = 5305.919
$ws.Range("I122").Value = 1253.4615
$ws.Range("K122").Value = 11281.1535
$ws.Range("M122").Value = -8831.153499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 14494.5
$ws.Range("I102").Value = 17659.334
$ws.Range("K102").Value = 17659.334
$ws.Range("M102").Value = -16037.334
$ws.Range("H132").Value = 3439.3257
$ws.Range("I132").Value = 3231.3235
$ws.Range("K132").Value = 9693.970499999999
$ws.Range("M132").Value = -7163.970499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20489.758
$ws.Range("J7").Value = 7689.7
$ws.Range("L7").Value = 7689.7
$ws.Range("N7").Value = -7913.7
$ws.Range("H16").Value = 15002.909
$ws.Range("J16").Value = 6001
$ws.Range("L16").Value = 6001
$ws.Range("N16").Value = -6341
$ws.Range("H22").Value = 782.375
$ws.Range("J22").Value = 781.1429000000001
$ws.Range("L22").Value = 781.1429000000001
$ws.Range("N22").Value = -1371.1429
$ws.Range("H27").Value = 782.375
$ws.Range("J27").Value = 781.1429000000001
$ws.Range("L27").Value = 781.1429000000001
$ws.Range("N27").Value = -995.1429000000001
$ws.Range("H40").Value = 21617.844
$ws.Range("I40").Value = 25002
$ws.Range("K40").Value = 25002
$ws.Range("M40").Value = -24866
$ws.Range("H46").Value = 1433.8636
$ws.Range("I46").Value = 807
$ws.Range("J46").Value = 2339.3333
$ws.Range("K46").Value = 807
$ws.Range("L46").Value = 2339.3333
$ws.Range("M46").Value = -619
$ws.Range("N46").Value = -2715.3333
$ws.Range("H126").Value = 20489.758
$ws.Range("J126").Value = 7689.7
$ws.Range("L126").Value = 23069.1
$ws.Range("N126").Value = -28009.1
$ws.Range("H132").Value = 423852.8
$ws.Range("I132").Value = 1032305.7
$ws.Range("K132").Value = 3096917.1
$ws.Range("M132").Value = -3094387.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 4497
$ws.Range("J55").Value = 4497
$ws.Range("L55").Value = 4497
$ws.Range("N55").Value = -5051
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 1999.4
$ws.Range("J113").Value = 20003
$ws.Range("K113").Value = 5998.200000000001
$ws.Range("L113").Value = 60009
$ws.Range("M113").Value = -3828.200000000001
$ws.Range("N113").Value = -64349
$ws.Range("H126").Value = 25093.055
$ws.Range("I126").Value = 28445.334
$ws.Range("J126").Value = 8331.666999999999
$ws.Range("K126").Value = 85336.00199999999
$ws.Range("L126").Value = 24995.001
$ws.Range("M126").Value = -82866.00199999999
$ws.Range("N126").Value = -29935.001
$ws.Range("H132").Value = 3708.7563
$ws.Range("I132").Value = 3680.224
$ws.Range("K132").Value = 11040.672
$ws.Range("M132").Value = -8510.672

Write-Host "Applied all updates"